{"js": "// Replace each arithmetic equation in the worksheet table with its updated\n// version, as described by the commit diff. The table is a 20-row x 5-column\n// grid where every cell holds exactly one equation (e.g. \"85+7=92\"). We walk\n// the cells in document order (row by row, left to right) and overwrite each\n// one's text with the corresponding replacement, which keeps this resilient\n// even though most of the new values do not numerically match the old ones.\n\nconst replacements = [\n    [\"85+7=92\", \"41+18=59\"],\n    [\"50+4=54\", \"55+33=88\"],\n    [\"18+34=52\", \"12+10=22\"],\n    [\"5+52=57\", \"56+37=93\"],\n    [\"1+93=94\", \"12+60=72\"],\n    [\"22+75=97\", \"22+45=67\"],\n    [\"72+5=77\", \"79-65=14\"],\n    [\"94-21=73\", \"66-49=17\"],\n    [\"37+25=62\", \"67+23=90\"],\n    [\"8+4=12\", \"22+69=91\"],\n    [\"27+17=44\", \"89-52=37\"],\n    [\"16+28=44\", \"75-68=7\"],\n    [\"38-17=21\", \"15+4=19\"],\n    [\"25+9=34\", \"7+18=25\"],\n    [\"55+36=91\", \"88-74=14\"],\n    [\"73-7=66\", \"11+60=71\"],\n    [\"19+63=82\", \"57-36=21\"],\n    [\"99-21=78\", \"33-30=3\"],\n    [\"97-79=18\", \"15+22=37\"],\n    [\"26+60=86\", \"2+94=96\"],\n    [\"93+5=98\", \"30+27=57\"],\n    [\"70-22=48\", \"56+40=96\"],\n    [\"5+8=13\", \"63+1=64\"],\n    [\"92-28=64\", \"82-57=25\"],\n    [\"49-25=24\", \"88-47=41\"],\n    [\"22-5=17\", \"66-19=47\"],\n    [\"90-17=73\", \"37+50=87\"],\n    [\"82-76=6\", \"39-3=36\"],\n    [\"48-36=12\", \"2+75=77\"],\n    [\"97-71=26\", \"11+4=15\"],\n    [\"34+39=73\", \"49-39=10\"],\n    [\"14-7=7\", \"88+1=89\"],\n    [\"12+29=41\", \"21+17=38\"],\n    [\"64+32=96\", \"61-47=14\"],\n    [\"39-15=24\", \"23-15=8\"],\n    [\"34-20=14\", \"70+28=98\"],\n    [\"17+36=53\", \"27-19=8\"],\n    [\"82-38=44\", \"74-36=38\"],\n    [\"41-27=14\", \"14+61=75\"],\n    [\"29-14=15\", \"95-47=48\"],\n    [\"64-18=46\", \"84-42=42\"],\n    [\"40+7=47\", \"73-14=59\"],\n    [\"98-78=20\", \"70-1=69\"],\n    [\"0+52=52\", \"97-86=11\"],\n    [\"81-0=81\", \"1+65=66\"],\n    [\"67-11=56\", \"96-30=66\"],\n    [\"15+24=39\", \"45-44=1\"],\n    [\"66-13=53\", \"50+38=88\"],\n    [\"80+3=83\", \"82-50=32\"],\n    [\"7+60=67\", \"68+15=83\"],\n    [\"12-1=11\", \"96-88=8\"],\n    [\"94-31=63\", \"79-20=59\"],\n    [\"67-56=11\", \"68-53=15\"],\n    [\"72+18=90\", \"62-7=55\"],\n    [\"47-29=18\", \"92-65=27\"],\n    [\"88-29=59\", \"76-65=11\"],\n    [\"57+32=89\", \"6-5=1\"],\n    [\"32+41=73\", \"85-1=84\"],\n    [\"42+26=68\", \"54-13=41\"],\n    [\"62-35=27\", \"18-17=1\"],\n    [\"29-22=7\", \"9+20=29\"],\n    [\"55-12=43\", \"36-24=12\"],\n    [\"58-28=30\", \"77-41=36\"],\n    [\"20-0=20\", \"12+21=33\"],\n    [\"24+2=26\", \"64-60=4\"],\n    [\"16+82=98\", \"70-13=57\"],\n    [\"27-24=3\", \"57-21=36\"],\n    [\"77-66=11\", \"6-0=6\"],\n    [\"6+75=81\", \"60-41=19\"],\n    [\"12+15=27\", \"93-44=49\"],\n    [\"32+30=62\", \"74-46=28\"],\n    [\"67-8=59\", \"67+26=93\"],\n    [\"48-0=48\", \"45+27=72\"],\n    [\"83-76=7\", \"61+38=99\"],\n    [\"92-11=81\", \"60-47=13\"],\n    [\"3+19=22\", \"54+5=59\"],\n    [\"14+46=60\", \"95-4=91\"],\n    [\"96-18=78\", \"8+23=31\"],\n    [\"26+11=37\", \"63-25=38\"],\n    [\"35+31=66\", \"53+1=54\"],\n    [\"11+87=98\", \"60-45=15\"],\n    [\"86-26=60\", \"0+40=40\"],\n    [\"32-7=25\", \"11+82=93\"],\n    [\"16-0=16\", \"59-24=35\"],\n    [\"45-41=4\", \"0+83=83\"],\n    [\"65+4=69\", \"97-36=61\"],\n    [\"4+71=75\", \"64-16=48\"],\n    [\"53+17=70\", \"17+70=87\"],\n    [\"71+10=81\", \"0+89=89\"],\n    [\"85-32=53\", \"84-7=77\"],\n    [\"32+36=68\", \"45-11=34\"],\n    [\"99-79=20\", \"98-33=65\"],\n    [\"50+21=71\", \"16-11=5\"],\n    [\"35+46=81\", \"90-60=30\"],\n    [\"11+54=65\", \"82+0=82\"],\n    [\"36+1=37\", \"18+70=88\"],\n    [\"38+28=66\", \"71+9=80\"],\n    [\"16+1=17\", \"84-69=15\"],\n    [\"96-94=2\", \"24+63=87\"],\n    [\"97-40=57\", \"0+51=51\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table of equations in the document body.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nconst cells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cells.push(cell);\n  }\n}\n\nif (cells.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} equation cells, found ${cells.length}.`\n  );\n}\n\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\n// Defensive check: confirm each cell still holds the equation we expect\n// before overwriting it, so a structural mismatch fails loudly instead of\n// silently clobbering the wrong cell.\ncells.forEach((cell, index) => {\n  const [oldText] = replacements[index];\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Cell ${index} expected \"${oldText}\" but found \"${cell.value}\".`\n    );\n  }\n});\n\ncells.forEach((cell, index) => {\n  const [, newText] = replacements[index];\n  cell.value = newText;\n});\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic equation in the worksheet table with its updated\n# version, as described by the commit diff. The document body is a single\n# table (20 rows x 5 columns) where every cell holds exactly one equation\n# (e.g. \"85+7=92\"). Each old equation string is unique in the document, so\n# we can safely drive the substitution with Find/Replace scoped to the whole\n# document content, one equation at a time.\n\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -lt 1) {\n    throw \"Expected a table of equations in the document body.\"\n}\n\n$replacements = @(\n    @{ Old = \"85+7=92\"; New = \"41+18=59\" },\n    @{ Old = \"50+4=54\"; New = \"55+33=88\" },\n    @{ Old = \"18+34=52\"; New = \"12+10=22\" },\n    @{ Old = \"5+52=57\"; New = \"56+37=93\" },\n    @{ Old = \"1+93=94\"; New = \"12+60=72\" },\n    @{ Old = \"22+75=97\"; New = \"22+45=67\" },\n    @{ Old = \"72+5=77\"; New = \"79-65=14\" },\n    @{ Old = \"94-21=73\"; New = \"66-49=17\" },\n    @{ Old = \"37+25=62\"; New = \"67+23=90\" },\n    @{ Old = \"8+4=12\"; New = \"22+69=91\" },\n    @{ Old = \"27+17=44\"; New = \"89-52=37\" },\n    @{ Old = \"16+28=44\"; New = \"75-68=7\" },\n    @{ Old = \"38-17=21\"; New = \"15+4=19\" },\n    @{ Old = \"25+9=34\"; New = \"7+18=25\" },\n    @{ Old = \"55+36=91\"; New = \"88-74=14\" },\n    @{ Old = \"73-7=66\"; New = \"11+60=71\" },\n    @{ Old = \"19+63=82\"; New = \"57-36=21\" },\n    @{ Old = \"99-21=78\"; New = \"33-30=3\" },\n    @{ Old = \"97-79=18\"; New = \"15+22=37\" },\n    @{ Old = \"26+60=86\"; New = \"2+94=96\" },\n    @{ Old = \"93+5=98\"; New = \"30+27=57\" },\n    @{ Old = \"70-22=48\"; New = \"56+40=96\" },\n    @{ Old = \"5+8=13\"; New = \"63+1=64\" },\n    @{ Old = \"92-28=64\"; New = \"82-57=25\" },\n    @{ Old = \"49-25=24\"; New = \"88-47=41\" },\n    @{ Old = \"22-5=17\"; New = \"66-19=47\" },\n    @{ Old = \"90-17=73\"; New = \"37+50=87\" },\n    @{ Old = \"82-76=6\"; New = \"39-3=36\" },\n    @{ Old = \"48-36=12\"; New = \"2+75=77\" },\n    @{ Old = \"97-71=26\"; New = \"11+4=15\" },\n    @{ Old = \"34+39=73\"; New = \"49-39=10\" },\n    @{ Old = \"14-7=7\"; New = \"88+1=89\" },\n    @{ Old = \"12+29=41\"; New = \"21+17=38\" },\n    @{ Old = \"64+32=96\"; New = \"61-47=14\" },\n    @{ Old = \"39-15=24\"; New = \"23-15=8\" },\n    @{ Old = \"34-20=14\"; New = \"70+28=98\" },\n    @{ Old = \"17+36=53\"; New = \"27-19=8\" },\n    @{ Old = \"82-38=44\"; New = \"74-36=38\" },\n    @{ Old = \"41-27=14\"; New = \"14+61=75\" },\n    @{ Old = \"29-14=15\"; New = \"95-47=48\" },\n    @{ Old = \"64-18=46\"; New = \"84-42=42\" },\n    @{ Old = \"40+7=47\"; New = \"73-14=59\" },\n    @{ Old = \"98-78=20\"; New = \"70-1=69\" },\n    @{ Old = \"0+52=52\"; New = \"97-86=11\" },\n    @{ Old = \"81-0=81\"; New = \"1+65=66\" },\n    @{ Old = \"67-11=56\"; New = \"96-30=66\" },\n    @{ Old = \"15+24=39\"; New = \"45-44=1\" },\n    @{ Old = \"66-13=53\"; New = \"50+38=88\" },\n    @{ Old = \"80+3=83\"; New = \"82-50=32\" },\n    @{ Old = \"7+60=67\"; New = \"68+15=83\" },\n    @{ Old = \"12-1=11\"; New = \"96-88=8\" },\n    @{ Old = \"94-31=63\"; New = \"79-20=59\" },\n    @{ Old = \"67-56=11\"; New = \"68-53=15\" },\n    @{ Old = \"72+18=90\"; New = \"62-7=55\" },\n    @{ Old = \"47-29=18\"; New = \"92-65=27\" },\n    @{ Old = \"88-29=59\"; New = \"76-65=11\" },\n    @{ Old = \"57+32=89\"; New = \"6-5=1\" },\n    @{ Old = \"32+41=73\"; New = \"85-1=84\" },\n    @{ Old = \"42+26=68\"; New = \"54-13=41\" },\n    @{ Old = \"62-35=27\"; New = \"18-17=1\" },\n    @{ Old = \"29-22=7\"; New = \"9+20=29\" },\n    @{ Old = \"55-12=43\"; New = \"36-24=12\" },\n    @{ Old = \"58-28=30\"; New = \"77-41=36\" },\n    @{ Old = \"20-0=20\"; New = \"12+21=33\" },\n    @{ Old = \"24+2=26\"; New = \"64-60=4\" },\n    @{ Old = \"16+82=98\"; New = \"70-13=57\" },\n    @{ Old = \"27-24=3\"; New = \"57-21=36\" },\n    @{ Old = \"77-66=11\"; New = \"6-0=6\" },\n    @{ Old = \"6+75=81\"; New = \"60-41=19\" },\n    @{ Old = \"12+15=27\"; New = \"93-44=49\" },\n    @{ Old = \"32+30=62\"; New = \"74-46=28\" },\n    @{ Old = \"67-8=59\"; New = \"67+26=93\" },\n    @{ Old = \"48-0=48\"; New = \"45+27=72\" },\n    @{ Old = \"83-76=7\"; New = \"61+38=99\" },\n    @{ Old = \"92-11=81\"; New = \"60-47=13\" },\n    @{ Old = \"3+19=22\"; New = \"54+5=59\" },\n    @{ Old = \"14+46=60\"; New = \"95-4=91\" },\n    @{ Old = \"96-18=78\"; New = \"8+23=31\" },\n    @{ Old = \"26+11=37\"; New = \"63-25=38\" },\n    @{ Old = \"35+31=66\"; New = \"53+1=54\" },\n    @{ Old = \"11+87=98\"; New = \"60-45=15\" },\n    @{ Old = \"86-26=60\"; New = \"0+40=40\" },\n    @{ Old = \"32-7=25\"; New = \"11+82=93\" },\n    @{ Old = \"16-0=16\"; New = \"59-24=35\" },\n    @{ Old = \"45-41=4\"; New = \"0+83=83\" },\n    @{ Old = \"65+4=69\"; New = \"97-36=61\" },\n    @{ Old = \"4+71=75\"; New = \"64-16=48\" },\n    @{ Old = \"53+17=70\"; New = \"17+70=87\" },\n    @{ Old = \"71+10=81\"; New = \"0+89=89\" },\n    @{ Old = \"85-32=53\"; New = \"84-7=77\" },\n    @{ Old = \"32+36=68\"; New = \"45-11=34\" },\n    @{ Old = \"99-79=20\"; New = \"98-33=65\" },\n    @{ Old = \"50+21=71\"; New = \"16-11=5\" },\n    @{ Old = \"35+46=81\"; New = \"90-60=30\" },\n    @{ Old = \"11+54=65\"; New = \"82+0=82\" },\n    @{ Old = \"36+1=37\"; New = \"18+70=88\" },\n    @{ Old = \"38+28=66\"; New = \"71+9=80\" },\n    @{ Old = \"16+1=17\"; New = \"84-69=15\" },\n    @{ Old = \"96-94=2\"; New = \"24+63=87\" },\n    @{ Old = \"97-40=57\"; New = \"0+51=51\" }\n)\n\n$replacedCount = 0\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        0,           # Wrap: wdFindStop\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace: wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Could not find expected equation text '$($pair.Old)' to replace.\"\n    }\n    $replacedCount++\n}\n\nif ($replacedCount -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) replacements, performed $replacedCount.\"\n}\n"}
